# Update "想去人数" (want-to-go count) figures in column F for both the
# "展览" and "全部类型" worksheets, for rows 3-11 and 17-18.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 269
    4  = 280
    5  = 823
    6  = 275
    7  = 6652
    8  = 55
    9  = 75
    10 = 118
    11 = 79
    17 = 558
    18 = 58
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
